# Add new data to the "DataSheet" sheet's deleteCustomer id list.
# The existing block (header "deleteCustomer"/"ID" + 21 customer ids,
# formerly at rows 24-46) is replaced by a new block starting 7 rows
# earlier (rows 17-18 header, rows 19-68 = 50 customer ids).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSheet")

$newIds = @(
    "cus_Kv3YAZpUxk7nBX","cus_Kv3Yj6sjAAuwsd","cus_Kv3Y0qTVOrYvhY","cus_Kv3YRJAn5ARQjp","cus_Kv3YdvFA4aOp5C",
    "cus_Kv3YKC4NITLVzw","cus_Kv3YJHH005LtWZ","cus_Kv3YScFptrL85n","cus_Kv3YgFXgOniq1Z","cus_Kv3YmBMJws74x4",
    "cus_Kv3YdhQ4XD17Nh","cus_Kv3YVf44eSRsVA","cus_Kv3YqxjpPoK7Bn","cus_Kv3YmyUsBPzHgQ","cus_Kv3YtMm6eQOin3",
    "cus_Kv3YhqUukhBoit","cus_Kv3Y4uMUZdppAl","cus_Kv3YUIBSeZHp10","cus_Kv3YSVvTkyREvk","cus_Kv3YbR2pbfKvab",
    "cus_Kv3YrFcLabBdar","cus_Kv3YWLCtW85X2Y","cus_Kv3YrEzvy3YgRP","cus_Kv3Y7JwQ4LMALT","cus_Kv3YAk4S69Ifee",
    "cus_Kv3YqNAL7zSrkQ","cus_Kv3YWf5ahkCrGa","cus_Kv3YuD89KMvKWz","cus_Kv3YQYHCxSS5hs","cus_Kv3Y9qZyBll4Kn",
    "cus_Kv3YYOW0c9t5p4","cus_Kv3YhYUEJinHoD","cus_Kv3XyI4k0v8uUn","cus_Kv3XnrzOJfZljq","cus_Kv3XDmjIyBmMjO",
    "cus_Kv3XbB5jhWGe2M","cus_Kv3X60zEj4MBwp","cus_Kv3X39NWA8Ubfk","cus_Kv3Xz8RBr65Yrx","cus_Kv3XBze2VYr2nE",
    "cus_Kv3XTji6ESoO4d","cus_Kv3XDAdGfeyeUR","cus_Kv3XgsO1MqgBuF","cus_Kv3XGlgLsZUuTV","cus_Kv3Xj195ajxsyY",
    "cus_Kv3XhPav2vvMir","cus_Kv3XaRpjiTzUEL","cus_Kv3Xjv27Ut5Y1L","cus_Kv3Xckyn3mTVwf","cus_Kv3X8vwpwm8OdV"
)

# Template cell that already carries the bordered look used by the
# "deleteCustomer" header and by every id row (A18 - the "ID" sub-header -
# intentionally keeps the plain/unbordered look, so it is left untouched).
$ws.Range("A24").Copy()
$ws.Range("A17").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("A19:A68").PasteSpecial(-4122)   # xlPasteFormats

# Header rows (shifted up from 24/25 to 17/18)
$ws.Range("A17").Value = "deleteCustomer"
$ws.Range("A18").Value = "ID"

# 50 new customer ids replacing the old 21
for ($i = 0; $i -lt $newIds.Length; $i++) {
    $ws.Cells.Item(19 + $i, 1).Value = $newIds[$i]
}

$ws.Activate()
$ws.Range("A18").Select()
